$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item('展览')
$ws1.Range('F2').Value = 2811
$ws1.Range('G2').Value = 70
$ws1.Range('F3').Value = 345
$ws1.Range('F4').Value = 366
$ws1.Range('G4').Value = 98
$ws1.Range('G5').Value = 88
$ws1.Range('G6').Value = 70
$ws1.Range('F7').Value = 328
$ws1.Range('G7').Value = '不可售'
$ws1.Range('F8').Value = 549
$ws1.Range('C9').Value = '杭州·第十届次元鹿角动漫游戏展（取消）'
$ws1.Range('D9').Value = '万融城3幢1楼 头号玩家数字运动俱乐部'
$ws1.Range('E9').Value = '2024.03.16 10:00-03.17 17:00'
$ws1.Range('F9').Value = 1161
$ws1.Range('G9').Value = '不可售'
$ws1.Range('H9').Value = 'https://show.bilibili.com/platform/detail.html?id=81136'
$ws1.Range('I9').Value = '//i1.hdslb.com/bfs/openplatform/202401/w8iHjfOA1705651976885.jpeg'
$ws1.Range('F10').Value = 554
$ws1.Range('F11').Value = 9500
$ws1.Range('F12').Value = 408
$ws1.Range('F13').Value = 2510
$ws1.Range('F17').Value = 473
$ws1.Range('F18').Value = 672
$ws1.Range('F21').Value = 1005
$ws1.Range('F22').Value = 2962
$ws1.Range('F23').Value = 2237
$ws1.Range('F25').Value = 1928
$ws1.Range('F26').Value = 1935
$ws1.Range('F27').Value = 484
$ws1.Range('F31').Value = 175
$ws1.Range('F34').Value = 340
$ws1.Range('F36').Value = 307
$ws1.Range('F39').Value = 115
$ws1.Range('F40').Value = 1421
$ws1.Range('F41').Value = 121
$ws1.Range('F42').Value = 1470
$ws1.Range('F44').Value = 339
$ws1.Range('F46').Value = 361
$ws1.Range('F47').Value = 733

$ws2 = $wb.Worksheets.Item('演出')
$ws2.Range('F3').Value = 24

$ws4 = $wb.Worksheets.Item('全部类型')
$ws4.Range('F2').Value = 2811
$ws4.Range('G2').Value = 70
$ws4.Range('F3').Value = 366
$ws4.Range('G3').Value = 98
$ws4.Range('G4').Value = 88
$ws4.Range('G6').Value = 70
$ws4.Range('F7').Value = 549
$ws4.Range('B8').Value = '2024-03-23'
$ws4.Range('C8').Value = '杭州·AD02动漫展'
$ws4.Range('D8').Value = '浙江省杭州市萧山区奔竞大道353号 国际博览中心'
$ws4.Range('E8').Value = '2024.03.23 10:00-03.24 17:00'
$ws4.Range('F8').Value = 9500
$ws4.Range('G8').Value = 75
$ws4.Range('H8').Value = 'https://show.bilibili.com/platform/detail.html?id=80905'
$ws4.Range('I8').Value = '//i1.hdslb.com/bfs/openplatform/202401/D3QaPamg1705397424553.jpeg'
$ws4.Range('C9').Value = '杭州·AD02动漫展  青柳尊哉内场票'
$ws4.Range('E9').Value = '2024.03.23 10:00-03.23 17:00'
$ws4.Range('F9').Value = 408
$ws4.Range('G9').Value = 528
$ws4.Range('H9').Value = 'https://show.bilibili.com/platform/detail.html?id=81503'
$ws4.Range('I9').Value = '//i1.hdslb.com/bfs/openplatform/202401/OmqxboDC1706522627528.jpeg'
$ws4.Range('C10').Value = '杭州·AD02动漫展--卡琳娜签售票'
$ws4.Range('E10').Value = '2024.03.23 09:30-03.23 17:00'
$ws4.Range('F10').Value = 2510
$ws4.Range('G10').Value = '已售罄'
$ws4.Range('H10').Value = 'https://show.bilibili.com/platform/detail.html?id=81941'
$ws4.Range('I10').Value = '//i1.hdslb.com/bfs/openplatform/202402/CZjxY9ZC1708416661613.jpeg'
$ws4.Range('F11').Value = 24
$ws4.Range('F16').Value = 672
$ws4.Range('F18').Value = 1005
$ws4.Range('F19').Value = 2962
$ws4.Range('F20').Value = 2237
$ws4.Range('F21').Value = 1928
$ws4.Range('F22').Value = 484
$ws4.Range('F26').Value = 175
$ws4.Range('F29').Value = 340
$ws4.Range('F31').Value = 307
$ws4.Range('F37').Value = 115
$ws4.Range('F38').Value = 1421
$ws4.Range('F40').Value = 121
$ws4.Range('F41').Value = 1470
$ws4.Range('F44').Value = 339
$ws4.Range('F46').Value = 361
$ws4.Range('F47').Value = 733
